$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New bill 13 (rows 17-18): Transfer Mandiri, 26 Oct 2024 17:32:19 ---
$ws.Range("A17").Value = 13
$ws.Range("B17").Value = "sudah"
$ws.Range("C17").Value = "26 October 2024 17:32:19"
$ws.Range("D17").Value = "Transfer Mandiri"
$ws.Range("E17").Value = 40000
$ws.Range("F17").Value = 40000
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 19
$ws.Range("I17").Value = 13
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = "Nasi Ayam Nashville BBQ"
$ws.Range("L17").Value = 20000
$ws.Range("M17").Value = 1
$ws.Range("N17").Value = 20000

$ws.Range("H18").Value = 20
$ws.Range("I18").Value = 13
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = "Nasi Ayam Nashville Sambal Bawang"
$ws.Range("L18").Value = 20000
$ws.Range("M18").Value = 1
$ws.Range("N18").Value = 20000

# --- New bill 14 (row 20): Transfer BCA, 26 Oct 2024 17:36:27 ---
$ws.Range("A20").Value = 14
$ws.Range("B20").Value = "sudah"
$ws.Range("C20").Value = "26 October 2024 17:36:27"
$ws.Range("D20").Value = "Transfer BCA"
$ws.Range("E20").Value = 20000
$ws.Range("F20").Value = 20000
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 21
$ws.Range("I20").Value = 14
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = "Mie Ayam Jakarta"
$ws.Range("L20").Value = 20000
$ws.Range("M20").Value = 1
$ws.Range("N20").Value = 20000

# --- New bill 15 (rows 22-23): Transfer Mandiri, 27 Oct 2024 14:37:27 ---
$ws.Range("A22").Value = 15
$ws.Range("B22").Value = "sudah"
$ws.Range("C22").Value = "27 October 2024 14:37:27"
$ws.Range("D22").Value = "Transfer Mandiri"
$ws.Range("E22").Value = 37000
$ws.Range("F22").Value = 37000
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 22
$ws.Range("I22").Value = 15
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = "Nasi Ayam Nashville Sambal Bawang"
$ws.Range("L22").Value = 20000
$ws.Range("M22").Value = 1
$ws.Range("N22").Value = 20000

$ws.Range("H23").Value = 23
$ws.Range("I23").Value = 15
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = "Sandwich Ayam Nashville Sambal Matah"
$ws.Range("L23").Value = 17000
$ws.Range("M23").Value = 1
$ws.Range("N23").Value = 17000

# --- Grand total block moved from row 18 down to row 26 ---
$ws.Range("J26").Value = "Total"
$ws.Range("N26").Value = 219000

# --- Subtotal-per-Jenis-Pembayaran block moved from rows 20-26 down to rows 28-34 ---
$ws.Range("A28").Value = "Total dari masing - masing Jenis Pembayaran:"

$ws.Range("A29").Value = "Cash"
$ws.Range("B29").Value = 0

$ws.Range("A30").Value = "Transfer Mandiri"
$ws.Range("B30").Value = 134000

$ws.Range("A31").Value = "Transfer BCA"
$ws.Range("B31").Value = 40000

$ws.Range("A32").Value = "QRIS"
$ws.Range("B32").Value = 45000

$ws.Range("A33").Value = "OVO"
$ws.Range("B33").Value = 0

$ws.Range("A34").Value = "Gopay"
$ws.Range("B34").Value = 0

# --- Clear the stale tail of the old subtotal block that used to occupy rows 21-26 ---
$ws.Range("A21").ClearContents()
$ws.Range("B21").ClearContents()
$ws.Range("A23").ClearContents()
$ws.Range("B23").ClearContents()
$ws.Range("A24").ClearContents()
$ws.Range("B24").ClearContents()
$ws.Range("A25").ClearContents()
$ws.Range("B25").ClearContents()
$ws.Range("A26").ClearContents()
$ws.Range("B26").ClearContents()
